# OpenTbs 1.8.1 beta - new common keywords for changing cell types (xlsxNum/xlsxBool/xlsxDate -> tbs:num/tbs:bool/tbs:date)
# plus a new named cell on the "Delete me" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Examples part 1": rename the "xlsx*" keywords to "tbs:*" ---
$ws1 = $wb.Worksheets.Item("Examples part 1")

$ws1.Range("C26").Value = "[cell2.score;block=tbs:cell;ope=tbs:num]"

$ws1.Range("C34").Value = "tbs:num"
$ws1.Range("C35").Value = "tbs:bool"
$ws1.Range("C36").Value = "tbs:date"

$ws1.Range("D36").Value = "[onshow.x_dt;ope=tbs:date]"
$ws1.Range("D35").Value = "[onshow.x_bt;ope=tbs:bool]"
$ws1.Range("D34").Value = "[onshow.x_num;ope=tbs:num]"

$ws1.Range("E20").Value = "[a.score;ope=tbs:num]"
$ws1.Range("F20").Value = "[a.score;ope=tbs:num]"

# --- Sheet "Delete me": add a new named cell B6 ---
$wsDel = $wb.Worksheets.Item("Delete me")
$wsDel.Range("B6").Value = "And this named cell too."

# F19 header changes from "Score" to "Score again"
$ws1.Range("F19").Value = "Score again"

# --- Workbook level: define the named range pointing at the new cell ---
$wb.Names.Add("the_named_cell", "='Delete me'!`$B`$6")
